$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamps (column A) for rows 7-13 and add new row 13
$timestamps = @{
    7 = "2024-08-30 16:01:26"
    8 = "2024-08-30 16:02:53"
    9 = "2024-08-30 16:13:05"
    10 = "2024-08-30 16:18:04"
    11 = "2024-08-30 16:19:59"
    12 = "2024-08-30 16:22:16"
    13 = "2024-08-30 16:53:33"
}

# Data values for columns B(2) through AA(27) for rows 7-13
$rowData = @{
    7 = @(15.38461538461539,0,0,0,7.692307692307693,0,25,0,30,0,0,0,0,10,11.11111111111111,0,0,25,25,7.692307692307693,100,100,0,0,0,0)
    8 = @(66.66666666666666,75,100,41.17647058823529,49.0566037735849,75,80,55.00000000000001,75,0,42.85714285714285,47.82608695652174,0,51.61290322580645,54.16666666666666,66.66666666666666,0,50,57.14285714285714,50.90909090909091,40,100,83.33333333333334,0,45.45454545454545,0)
    9 = @(70,60,33.33333333333333,75,48.27586206896552,100,33.33333333333333,66.66666666666666,68.42105263157895,0,0,53.84615384615385,66.66666666666666,50,60,80,0,50,80,68.75,85.71428571428571,100,100,0,75,0)
    10 = @(0,0,0,0,20,0,0,100,60,0,0,50,0,0,0,0,0,0,100,100,0,0,0,0,100,0)
    11 = @(0,0,0,100,100,0,0,0,0,0,0,33.33333333333333,0,0,100,0,0,0,0,100,100,0,100,0,0,0)
    12 = @(0,0,0,100,0,0,0,0,0,0,0,33.33333333333333,0,0,50,0,0,0,0,100,100,0,100,0,100,0)
    13 = @(0,0,0,0,33.33333333333333,0,0,0,0,0,0,50,0,50,0,100,0,0,0,0,0,0,100,0,0,100)
}

foreach ($row in 7..13) {
    $ws.Cells.Item($row, 1).Value = $timestamps[$row]
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}
